$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New ".58.com" subdomains appended under the existing "58同城" (58.com) entries.
$domains = @(
    "hz.58.com",
    "tj.58.com",
    "infotopweb.union.vip.58.com",
    "su.58.com",
    "lieche.58.com",
    "nj.58.com",
    "dg.58.com",
    "dl.58.com",
    "hrb.58.com",
    "chefenqi.58.com",
    "sjz.58.com",
    "weizhang.58.com",
    "hshi.58.com",
    "caipiao.58.com",
    "epost.58.com",
    "xa.58.com",
    "cha.58.com",
    "cq.58.com",
    "ez.58.com",
    "sou.58.com",
    "cs.58.com",
    "jn.58.com"
)

$company = $ws.Cells.Item(804, 2).Value()

$row = 805
foreach ($d in $domains) {
    $ws.Cells.Item($row, 1).Value = $d
    $ws.Cells.Item($row, 2).Value = $company
    $row = $row + 1
}

# Reset the view: scroll back to the top and select B18 (matches the saved view state).
$ws.Activate()
$ws.Range("B18").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 1
